# ISS-46: Fix reading group TextFrame
#
# Adds a new group shape ("Группа 3") containing two text boxes
# ("TextBox 2" / "ch1" and "TextBox 4" / "ch2") to slide 1, mirroring
# the think-cell style caption group added in the reference commit.

$EMU = 12700
# The COM layer here rounds Left/Top/Width/Height (which are expressed
# in points) down to the nearest EMU when they go through a property
# *setter* after shape creation, so a tiny epsilon nudges the value back
# onto the correct integral EMU boundary without visibly changing it.
$EPS = 0.00003
function Pt([double]$emuValue) {
    return ($emuValue / $EMU) + $EPS
}

$ppt2 = $ppt
$p = $ppt2.ActivePresentation
$s = $p.Slides.Item(1)

# --- child textbox 1 ("ch1") -------------------------------------------------
$tb1 = $s.Shapes.AddTextbox(1, (Pt 3779912), (Pt 908720), (Pt 1080120), (Pt 360040))
$tb1.Name = "TextBox 2"

$tf1 = $tb1.TextFrame
$tf1.WordWrap = -1
$tf1.AutoSize = 0
$tf1.MarginLeft = 0
$tf1.MarginRight = 0
$tf1.MarginTop = 0
$tf1.MarginBottom = 0

$tr1 = $tf1.TextRange
$tr1.Text = "ch1"
$tr1.Font.Name = "Georgia"
$tr1.Font.Size = 20
$tr1.ParagraphFormat.SpaceAfter = 9
$tf1.Ruler.Levels.Item(1).FirstMargin = -21.6

# --- child textbox 2 ("ch2") -------------------------------------------------
$tb2 = $s.Shapes.AddTextbox(1, (Pt 3779912), (Pt 1484784), (Pt 1080120), (Pt 360040))
$tb2.Name = "TextBox 4"

$tf2 = $tb2.TextFrame
$tf2.WordWrap = -1
$tf2.AutoSize = 0
$tf2.MarginLeft = 0
$tf2.MarginRight = 0
$tf2.MarginTop = 0
$tf2.MarginBottom = 0

$tr2 = $tf2.TextRange
$tr2.Text = "ch2"
$tr2.Font.Name = "Georgia"
$tr2.Font.Size = 20
$tr2.ParagraphFormat.SpaceAfter = 9
$tf2.Ruler.Levels.Item(1).FirstMargin = -21.6

# --- group the two textboxes and move the group into place -------------------
$range = $s.Shapes.Range(@($tb1.Name, $tb2.Name))
$grp = $range.Group()
$grp.Name = "Группа 3"
$grp.Left = Pt 2267744
$grp.Top = Pt 476672
